$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = "Jesse Hare"
$ws.Range("E1").Value = 8

$ws.Range("A3").Value = "Project Build"
$ws.Range("B3").Value = "Implementation of search functionality"
$ws.Range("C3").Value = 30
$ws.Range("D3").Value = 20

$ws.Range("A14").Value = "Cumulative Total:160"

$ws.Range("C10").Select()
